## fix issues de compilation 2cfb2c2de19ac0faf3624fc5289f4e36a65be006
##
## Metadata sheet: insert a new "Jurisdiction" property row right before
## "Description", and refresh the "Date" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Description" currently lives on row 11; insert a blank row above it so
# everything from "Description" down shifts from row N to row N+1.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the rest of the property/value rows (copy the
# style from the row immediately below, now holding "Description").
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh publish date.
$ws.Range("B8").Value = "2024-09-12T14:01:50+00:00"
